$d = $word.ActiveDocument

# The second table holds the test-case rows; column 6 ("Performed by/ Date")
# of each data row (rows 2-7) currently contains a single paragraph with
# just the date "3/5/16". Insert a new paragraph "Edgar Dimanarig" right
# before that date paragraph in every one of those rows.
$t2 = $d.Tables.Item(2)
for ($r = 2; $r -le 7; $r++) {
    $cell = $t2.Rows.Item($r).Cells.Item(6)
    $datePara = $cell.Range.Paragraphs.Item(1)
    $insertPoint = $datePara.Range.Duplicate
    $insertPoint.Collapse(1)
    $insertPoint.InsertBefore("Edgar Dimanarig" + [char]13)
}

# Fix the typo in the description of test case 15-6.
$d.Content.Find.Execute("The project start after the End date", $true, $false, $false, $false, $false, $true, 1, $false, "The project starts after the End date", 2)
